$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.483.45'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.911.83'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.96'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +9.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.40'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.32%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.717'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000336'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.11'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.536.09'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.26'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.913.29'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +9.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.02'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.396.53'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '428.66'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '88.43'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.15'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.02'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +9.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.48'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.63'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.40'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '684.12'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.17'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.85'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '68.88'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +13.38%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +11.95%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0863'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.89%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.40%  '
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.17'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.64%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.19'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.15'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.78%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.141'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.34'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0363'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +14.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.00'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +6.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.755.20'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +13.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.01'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.21%  '
